# The price list sheet ("Hoja1") is the active sheet in this workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the printed date in A1 (2024-05-03 -> 2024-05-24).
$ws.Range("A1").Value = 45436

# Refresh the unit prices in column D for the listed items (rows 34-37).
$ws.Range("D34").Value = 206.846
$ws.Range("D35").Value = 293.075
$ws.Range("D36").Value = 396.525
$ws.Range("D37").Value = 653.796
